{"js": "// Thank You Template \u2014 pathing/resources update:\n//  1. Insert two new blank paragraphs at the very top of the document\n//     (before the letterhead picture paragraph).\n//  2. Remove three of the four blank paragraphs that sit between the\n//     picture paragraph and the DATE field paragraph (leave one).\n//  3. Refresh the cached DATE field result from \"October 5, 2021\" to\n//     \"October 9, 2021\".\n//  4. Collapse the three runs that spell out the `r fu` template tag\n//     into a single run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Two new empty paragraphs before the existing first paragraph ---\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertParagraph(\"\", \"Before\");\nfirstParagraph.insertParagraph(\"\", \"Before\");\nawait context.sync();\n\n// --- 2. Trim the run of blank paragraphs ahead of the date paragraph ---\n// After the two inserts above, the body now reads:\n//   [0] new blank, [1] new blank, [2] picture, [3..6] four blanks, [7] date\n// Keep paragraph [3] and drop the next three blanks.\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[4].delete();\nparagraphs.items[4].delete();\nparagraphs.items[4].delete();\nawait context.sync();\n\n// --- 3. Update the stale DATE field text ---\nconst dateHits = body.search(\"October 5, 2021\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"October 9, 2021\", \"Replace\");\n  await context.sync();\n}\n\n// --- 4. Merge the split \"`r fu`\" template-tag runs into one run ---\nconst fuHits = body.search(\"`r fu`\", { matchCase: true });\nfuHits.load(\"items\");\nawait context.sync();\nif (fuHits.items.length > 0) {\n  fuHits.items[0].insertText(\"`r fu`\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Thank You Template \u2014 pathing/resources update:\n#  1. Insert two new blank paragraphs at the very top of the document\n#     (before the letterhead picture paragraph).\n#  2. Remove three of the four blank paragraphs that sit between the\n#     picture paragraph and the DATE field paragraph (leave one).\n#  3. Refresh the cached DATE field result from \"October 5, 2021\" to\n#     \"October 9, 2021\".\n#  4. Collapse the three runs that spell out the `r fu` template tag\n#     into a single run.\n\n$d = $word.ActiveDocument\n\n# --- 1. Two new empty paragraphs before the existing first paragraph ---\n$d.Paragraphs.Item(1).Range.InsertParagraphBefore()\n$d.Paragraphs.Item(1).Range.InsertParagraphBefore()\n\n# --- 2. Trim the run of blank paragraphs ahead of the date paragraph ---\n# Layout is now: [1] blank, [2] blank, [3] picture, [4..7] four blanks, [8] date.\n# Keep paragraph 4 and drop the next three blanks.\n$d.Paragraphs.Item(5).Range.Delete()\n$d.Paragraphs.Item(5).Range.Delete()\n$d.Paragraphs.Item(5).Range.Delete()\n\n# --- 3. Update the stale DATE field text ---\n$find = $d.Content.Find\n$find.Execute(\n    \"October 5, 2021\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"October 9, 2021\", 2\n)\n\n# --- 4. Merge the split \"`r fu`\" template-tag runs into one run ---\n$find2 = $d.Content.Find\n$find2.Execute(\n    \"``r fu``\", $false, $false, $false, $false, $false, $true, 1, $false,\n    \"``r fu``\", 2\n)\n"}
